# Added TestCase in Card Unit Testing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The existing "Last Runtime" dates in E2:E6 move from 8/5/2021 (44413)
#    to 8/6/2021 (44414).
$ws.Range("E2:E6").Value = 44414

# 2. Widen column B slightly (52.8984375 -> ~63.6 chars).
$ws.Columns.Item(2).ColumnWidth = 62.86

# 3. Append the new CARD_006 test case as row 7. Copy the date cell's
#    number-format/style from E6 first so the new E7 matches the existing
#    "m/d/yyyy" styled cells instead of creating a brand-new style entry.
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("A7").Value = "CARD_006"
$ws.Range("B7").Value = "Doesn't Change Card Data In Memory if didn’t receive ADMIN from Terminal"
$ws.Range("C7").Value = "Card Data  Isn't Changed in memory"
$ws.Range("D7").Value = "PASSED"
$ws.Range("E7").Value = 44414

# 4. Grow Table1 so the new row is included in the table range.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E7"))

# 5. Selection moves to E2.
$null = $ws.Range("E2").Select()
